# Applies the "Add files via upload" edit to CasosColombia.xlsx:
#  - A handful of existing cells in the data table had their value replaced
#    by the text "NaN" (or, in two cases, a "NaN" cell was replaced with a
#    real number), matching the existing shared string used elsewhere in
#    the sheet.
#  - A brand-new row (row 189, dated 2020-09-09 / serial 44083) was
#    appended with a full set of per-department case counts. One of its
#    cells (AM189) holds the literal text "####" which becomes a new
#    shared-string table entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell corrections in the existing data -----------------
$ws.Range("DD25").Value  = "NaN"
$ws.Range("CC44").Value  = "NaN"
$ws.Range("O51").Value   = 7
$ws.Range("X73").Value   = "NaN"
$ws.Range("AP84").Value  = "NaN"
$ws.Range("CY85").Value  = "NaN"
$ws.Range("O87").Value   = "NaN"
$ws.Range("AF127").Value = 6
$ws.Range("AF128").Value = "NaN"
$ws.Range("CF144").Value = "NaN"

# --- New row 189 --------------------------------------------------------
$ws.Range("A189").Value = 44083
$ws.Range("B189").Value = 987851
$ws.Range("C189").Value = 5230
$ws.Range("D189").Value = 126019
$ws.Range("E189").Value = 118123
$ws.Range("F189").Value = 334403
$ws.Range("G189").Value = 44587
$ws.Range("H189").Value = 5698
$ws.Range("I189").Value = 4609
$ws.Range("J189").Value = 8185
$ws.Range("K189").Value = 8821
$ws.Range("L189").Value = 17094
$ws.Range("M189").Value = 6852
$ws.Range("N189").Value = 27262
$ws.Range("O189").Value = 35611
$ws.Range("P189").Value = 8356
$ws.Range("Q189").Value = 7663
$ws.Range("R189").Value = 19767
$ws.Range("S189").Value = 13412
$ws.Range("T189").Value = 23616
$ws.Range("U189").Value = 15593
$ws.Range("V189").Value = 4237
$ws.Range("W189").Value = 2080
$ws.Range("X189").Value = 9211
$ws.Range("Y189").Value = 26799
$ws.Range("Z189").Value = 18756
$ws.Range("AA189").Value = 11415
$ws.Range("AB189").Value = 76869
$ws.Range("AC189").Value = 1672
$ws.Range("AD189").Value = 455
$ws.Range("AE189").Value = 560
$ws.Range("AF189").Value = 456
$ws.Range("AG189").Value = 346
$ws.Range("AH189").Value = 220
$ws.Range("AI189").Value = 508
$ws.Range("AJ189").Value = 3630
$ws.Range("AK189").Value = 4905
$ws.Range("AL189").Value = 66491
$ws.Range("AM189").Value = "####"
$ws.Range("AN189").Value = 4599
$ws.Range("AO189").Value = 59907
$ws.Range("AP189").Value = 1605
$ws.Range("AQ189").Value = 36829
$ws.Range("AR189").Value = 2604
$ws.Range("AS189").Value = 10758
$ws.Range("AT189").Value = 2756
$ws.Range("AU189").Value = 2932
$ws.Range("AV189").Value = 6792
$ws.Range("AW189").Value = 2613
$ws.Range("AX189").Value = 1863
$ws.Range("AY189").Value = 4849
$ws.Range("AZ189").Value = 4942
$ws.Range("BA189").Value = 72011
$ws.Range("BB189").Value = 16862
$ws.Range("BC189").Value = 4334
$ws.Range("BD189").Value = 12108
$ws.Range("BE189").Value = 6027
$ws.Range("BF189").Value = 515
$ws.Range("BG189").Value = 2623
$ws.Range("BH189").Value = 4805
$ws.Range("BI189").Value = 1414
$ws.Range("BJ189").Value = 3413
$ws.Range("BK189").Value = 13103
$ws.Range("BL189").Value = 13540
$ws.Range("BM189").Value = 12826
$ws.Range("BN189").Value = 25874
$ws.Range("BO189").Value = 3649
$ws.Range("BP189").Value = 1464
$ws.Range("BQ189").Value = 12589
$ws.Range("BR189").Value = 10971
$ws.Range("BS189").Value = 11293
$ws.Range("BT189").Value = 2412
$ws.Range("BU189").Value = 2362
$ws.Range("BV189").Value = 4980
$ws.Range("BW189").Value = 4705
$ws.Range("BX189").Value = 1483
$ws.Range("BY189").Value = 6998
$ws.Range("BZ189").Value = 3917
$ws.Range("CA189").Value = 2014
$ws.Range("CB189").Value = 1193
$ws.Range("CC189").Value = 3138
$ws.Range("CD189").Value = 3184
$ws.Range("CE189").Value = 1827
$ws.Range("CF189").Value = 1537
$ws.Range("CG189").Value = 6756
$ws.Range("CH189").Value = 2032
$ws.Range("CI189").Value = 1774
$ws.Range("CJ189").Value = 1985
$ws.Range("CK189").Value = 2458
$ws.Range("CL189").Value = 2278
$ws.Range("CM189").Value = 2553
$ws.Range("CN189").Value = 1682
$ws.Range("CO189").Value = 1635
$ws.Range("CP189").Value = 1815
$ws.Range("CQ189").Value = 932
$ws.Range("CR189").Value = 4790
$ws.Range("CS189").Value = 1546
$ws.Range("CT189").Value = 1294
$ws.Range("CU189").Value = 1221
$ws.Range("CV189").Value = 1937
$ws.Range("CW189").Value = 1827
$ws.Range("CX189").Value = 1033
$ws.Range("CY189").Value = 1214
$ws.Range("CZ189").Value = 1347
$ws.Range("DA189").Value = 1878
$ws.Range("DB189").Value = 1468
$ws.Range("DC189").Value = 1714
$ws.Range("DD189").Value = 1372
$ws.Range("DE189").Value = 605
$ws.Range("DF189").Value = 579
$ws.Range("DG189").Value = 1041
$ws.Range("DH189").Value = 921
$ws.Range("DI189").Value = 682
$ws.Range("DJ189").Value = 809
$ws.Range("DK189").Value = 577
$ws.Range("DL189").Value = 921
$ws.Range("DM189").Value = 1297
$ws.Range("DN189").Value = 892
$ws.Range("DO189").Value = 791
$ws.Range("DP189").Value = 653
$ws.Range("DQ189").Value = 901
$ws.Range("DR189").Value = 207927
$ws.Range("DS189").Value = 409923
$ws.Range("DT189").Value = 16173
$ws.Range("DU189").Value = 172027
$ws.Range("DV189").Value = 116158
$ws.Range("DW189").Value = 42441
$ws.Range("DX189").Value = 12768

# Keep the active selection pointing at the new bottom-right cell, as in
# the authored workbook (selection moves to DY189 after the append).
$ws.Range("DY189").Select()
